$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2-10 (was rows 2-4)

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Serping1"
$ws.Cells.Item(2, 3).Value = "Selp"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 11.921572
$ws.Cells.Item(2, 8).Value = 35.764716
$ws.Cells.Item(2, 9).Value = 0.006971694289596158
$ws.Cells.Item(2, 10).Value = 0.006971694289596159
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 26.450162
$ws.Cells.Item(2, 14).Value = 79.350486
$ws.Cells.Item(2, 15).Value = 0.9821351879331711
$ws.Cells.Item(2, 16).Value = 0.9821351879331711
$ws.Cells.Item(2, 17).Value = 315.327510694664
$ws.Cells.Item(2, 18).Value = 2837.947596251976
$ws.Cells.Item(2, 19).Value = 0.006847146281325139
$ws.Cells.Item(2, 20).Value = 0.00684714628132514

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Serping1"
$ws.Cells.Item(3, 3).Value = "Selp"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 11.921572
$ws.Cells.Item(3, 8).Value = 35.764716
$ws.Cells.Item(3, 9).Value = 0.006971694289596158
$ws.Cells.Item(3, 10).Value = 0.006971694289596159
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.321934
$ws.Cells.Item(3, 14).Value = 0.965802
$ws.Cells.Item(3, 15).Value = 0.01195390446349922
$ws.Cells.Item(3, 16).Value = 0.01195390446349922
$ws.Cells.Item(3, 17).Value = 3.837959360248
$ws.Cells.Item(3, 18).Value = 34.541634242232
$ws.Cells.Item(3, 19).Value = 0.00008333896748655554
$ws.Cells.Item(3, 20).Value = 0.00008333896748655557

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Serping1"
$ws.Cells.Item(4, 3).Value = "Selp"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 11.921572
$ws.Cells.Item(4, 8).Value = 35.764716
$ws.Cells.Item(4, 9).Value = 0.006971694289596158
$ws.Cells.Item(4, 10).Value = 0.006971694289596159
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1591883333333333
$ws.Cells.Item(4, 14).Value = 0.477565
$ws.Cells.Item(4, 15).Value = 0.005910907603329674
$ws.Cells.Item(4, 16).Value = 0.005910907603329674
$ws.Cells.Item(4, 17).Value = 1.897775177393333
$ws.Cells.Item(4, 18).Value = 17.07997659654
$ws.Cells.Item(4, 19).Value = 0.00004120904078446401
$ws.Cells.Item(4, 20).Value = 0.00004120904078446401

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Serping1"
$ws.Cells.Item(5, 3).Value = "Selp"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1516.953124666667
$ws.Cells.Item(5, 8).Value = 4550.859374
$ws.Cells.Item(5, 9).Value = 0.8871089682487887
$ws.Cells.Item(5, 10).Value = 0.8871089682487888
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 26.450162
$ws.Cells.Item(5, 14).Value = 79.350486
$ws.Cells.Item(5, 15).Value = 0.9821351879331711
$ws.Cells.Item(5, 16).Value = 0.9821351879331711
$ws.Cells.Item(5, 17).Value = 40123.65589383953
$ws.Cells.Item(5, 18).Value = 361112.9030445557
$ws.Cells.Item(5, 19).Value = 0.8712609332482256
$ws.Cells.Item(5, 20).Value = 0.8712609332482257

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Serping1"
$ws.Cells.Item(6, 3).Value = "Selp"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1516.953124666667
$ws.Cells.Item(6, 8).Value = 4550.859374
$ws.Cells.Item(6, 9).Value = 0.8871089682487887
$ws.Cells.Item(6, 10).Value = 0.8871089682487888
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.321934
$ws.Cells.Item(6, 14).Value = 0.965802
$ws.Cells.Item(6, 15).Value = 0.01195390446349922
$ws.Cells.Item(6, 16).Value = 0.01195390446349922
$ws.Cells.Item(6, 17).Value = 488.3587872364386
$ws.Cells.Item(6, 18).Value = 4395.229085127948
$ws.Cells.Item(6, 19).Value = 0.01060441585515938
$ws.Cells.Item(6, 20).Value = 0.01060441585515939

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Serping1"
$ws.Cells.Item(7, 3).Value = "Selp"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1516.953124666667
$ws.Cells.Item(7, 8).Value = 4550.859374
$ws.Cells.Item(7, 9).Value = 0.8871089682487887
$ws.Cells.Item(7, 10).Value = 0.8871089682487888
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.1591883333333333
$ws.Cells.Item(7, 14).Value = 0.477565
$ws.Cells.Item(7, 15).Value = 0.005910907603329674
$ws.Cells.Item(7, 16).Value = 0.005910907603329674
$ws.Cells.Item(7, 17).Value = 241.4812396604789
$ws.Cells.Item(7, 18).Value = 2173.33115694431
$ws.Cells.Item(7, 19).Value = 0.005243619145403708
$ws.Cells.Item(7, 20).Value = 0.005243619145403709

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Serping1"
$ws.Cells.Item(8, 3).Value = "Selp"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 181.1216836666667
$ws.Cells.Item(8, 8).Value = 543.365051
$ws.Cells.Item(8, 9).Value = 0.1059193374616151
$ws.Cells.Item(8, 10).Value = 0.1059193374616151
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 26.450162
$ws.Cells.Item(8, 14).Value = 79.350486
$ws.Cells.Item(8, 15).Value = 0.9821351879331711
$ws.Cells.Item(8, 16).Value = 0.9821351879331711
$ws.Cells.Item(8, 17).Value = 4790.697874696088
$ws.Cells.Item(8, 18).Value = 43116.28087226479
$ws.Cells.Item(8, 19).Value = 0.1040271084036203
$ws.Cells.Item(8, 20).Value = 0.1040271084036203

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Serping1"
$ws.Cells.Item(9, 3).Value = "Selp"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 181.1216836666667
$ws.Cells.Item(9, 8).Value = 543.365051
$ws.Cells.Item(9, 9).Value = 0.1059193374616151
$ws.Cells.Item(9, 10).Value = 0.1059193374616151
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.321934
$ws.Cells.Item(9, 14).Value = 0.965802
$ws.Cells.Item(9, 15).Value = 0.01195390446349922
$ws.Cells.Item(9, 16).Value = 0.01195390446349922
$ws.Cells.Item(9, 17).Value = 58.30922810954466
$ws.Cells.Item(9, 18).Value = 524.783052985902
$ws.Cells.Item(9, 19).Value = 0.001266149640853281
$ws.Cells.Item(9, 20).Value = 0.001266149640853282

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Serping1"
$ws.Cells.Item(10, 3).Value = "Selp"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 181.1216836666667
$ws.Cells.Item(10, 8).Value = 543.365051
$ws.Cells.Item(10, 9).Value = 0.1059193374616151
$ws.Cells.Item(10, 10).Value = 0.1059193374616151
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1591883333333333
$ws.Cells.Item(10, 14).Value = 0.477565
$ws.Cells.Item(10, 15).Value = 0.005910907603329674
$ws.Cells.Item(10, 16).Value = 0.005910907603329674
$ws.Cells.Item(10, 17).Value = 28.83245895342389
$ws.Cells.Item(10, 18).Value = 259.492130580815
$ws.Cells.Item(10, 19).Value = 0.0006260794171415023
$ws.Cells.Item(10, 20).Value = 0.0006260794171415024
